$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New simulation run for GJ 9827 d: lower the latitude input (B5) from 40 to 20.
# K5's formula recalculates automatically.
$ws.Range("B5").Value = 20

# Record the resulting inclination angle as a standalone value snapshot in K10.
$ws.Range("K10").Value = 89.017641452899696

# Leave the active selection on the newly written cell, matching the saved view.
$ws.Range("K10").Select()
